$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "characteristicts" -> "characteristics"
#    The misspelling was wrapped in <w:proofErr spellStart/spellEnd>.
#    Replacing across the whole sentence collapses the three runs
#    (and drops the now-stale proofErr markers) into a single run.
# ---------------------------------------------------------------------
$sentence = $d.Content
$sentence.Find.Execute( `
    "one class takes over all characteristicts of another class", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "one class takes over all characteristics of another class", 2) | Out-Null

# Re-locate "characteristics" inside the just-fixed sentence (rather than
# a hard-coded offset) so we can split a run right after the word and
# park the moved _GoBack bookmark there.
$word1 = $sentence.Duplicate
$word1.Find.Execute("characteristics", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null

# Force a run boundary *before* "characteristics" with a throw-away
# bookmark, so the word ends up in its own <w:r> (matching the target
# layout), then remove the scratch bookmark once it has served its
# purpose.
$wordStart = $d.Range($word1.Start, $word1.Start)
$d.Bookmarks.Add("zzTmpSplit", $wordStart) | Out-Null

# Word keeps a single "last edit" _GoBack bookmark; adding it here moves
# it from wherever it used to be (the Miscellaneous Inheritance heading)
# to right after "characteristics".
$wordEnd = $d.Range($word1.End, $word1.End)
$d.Bookmarks.Add("_GoBack", $wordEnd) | Out-Null

$d.Bookmarks("zzTmpSplit").Delete()

# ---------------------------------------------------------------------
# 2) "Miscellaneous Inheritan" + "ce Situations" (previously split by
#    the old _GoBack bookmark) -> single run "Miscellaneous Inheritance
#    Situations" now that the bookmark moved away.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Miscellaneous Inheritance Situations", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Miscellaneous Inheritance Situations", 2) | Out-Null
